# Apply commit: "add task#3; comment on task#4"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task #4 breakdown-of-work comment update ---
# B93: "If you worked with a partner, describe the breakdown of work"
$ws.Range("C93").Value = 'Bo made an improvement in this part; Pokai helped with coding'

# --- Task #3 (Adadelta) section: answer the two previously-blank questions ---
# B47: "What made you believe it would improve the model?  "
# B48: "In what way(s) did you expect it to impact the model?"
$ws.Range("C47").Value = 'For this, I believe we have to go through Adagrad first since Adadelta is actually an extension of Adagrad. Adadelta, if compared with GD, essentially associates learning rates'' change with features'' frequencies. That is, Adadelta would update with larger learning rate for parameters associated with infrequent features; on the other hand, it would update with smaller learning rate for those associated with frequent features. And thus, this algorithm might be more suitable when data sparseness increases. If one considers the size of training data being provided (2487) and the number of pixels for each image on three channels (320 X 240 X 3), one would tend to conjecture that there would be some data sparseness issues in this training task. Adadelta is a variant of Adagrad that essentially imposes a regularization coefficient in searching of the parameters.'
$ws.Range("C48").Value = 'For this, I believe we have to go through Adagrad first since Adadelta is actually an extension of Adagrad. Adadelta, if compared with GD, essentially associates learning rates'' change with features'' frequencies. That is, Adadelta would update with larger learning rate for parameters associated with infrequent features; on the other hand, it would update with smaller learning rate for those associated with frequent features. And thus, this algorithm might be more suitable when data sparseness increases. If one considers the size of training data being provided (2487) and the number of pixels for each image on three channels (320 X 240 X 3), one would tend to conjecture that there would be some data sparseness issues in this training task. Adadelta is a variant of Adagrad that essentially imposes a regularization coefficient in searching of the parameters.'

# --- Restore the view/selection state as last saved ---
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D49").Select() | Out-Null
